$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Market share_class_min ZEV values for 2035-2050 from 0.99 to 1
$ws.Range("T3:W3").Value = 1

# Update the sheet view selection range to match the actual used range A1:X7
$ws.Range("A1:X7").Select()
